$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (Petition for GCS)
$ws.Range("B2").Value = 0.694189602446483
$ws.Range("C2").Value = 0.670241286863271
$ws.Range("D2").Value = 0.776119402985075
$ws.Range("E2").Value = 0.70404984423676
$ws.Range("F2").Value = 0.50297176820208

# Update row 3 values (Petition for NR)
$ws.Range("B3").Value = 0.645051194539249
$ws.Range("C3").Value = 0.65625
$ws.Range("D3").Value = 0.752727272727273
$ws.Range("E3").Value = 0.656346749226006
$ws.Range("F3").Value = 0.547417116422513

# Row 4 becomes "Petition (any)" with new values, old rows 4/5/6 data replaced
$ws.Range("A4").Value = "Petition (any)"
$ws.Range("B4").Value = 0.670967741935484
$ws.Range("C4").Value = 0.663143989431968
$ws.Range("D4").Value = 0.76427255985267
$ws.Range("E4").Value = 0.680124223602484
$ws.Range("F4").Value = 0.524782444192206

# Remove old rows 5 and 6 (order_petition_271 row and old Petition (any) row)
$ws.Range("A5:F6").Delete()
